$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = "dsada"
$ws.Range("D2").Value = "dsada"
$ws.Range("E2").Value = 333
$ws.Range("F2").Value = "0:0"
$ws.Range("G2").Value = "Плацкарт"

$ws.Range("C3").Value = "dsada"
$ws.Range("D3").Value = "sdasdsada"
$ws.Range("E3").Value = 2222
$ws.Range("F3").Value = "0:0"
$ws.Range("G3").Value = "Плацкарт"

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "dsada"
$ws.Range("D4").Value = "sdasdsada"
$ws.Range("E4").Value = 2221
$ws.Range("F4").Value = "0:0"
$ws.Range("G4").Value = "Плацкарт"

# Row 5 updates
$ws.Range("C5").Value = "dsadasd"
$ws.Range("D5").Value = "sadasda"
$ws.Range("E5").Value = 1000
$ws.Range("F5").Value = "0:0"
$ws.Range("G5").Value = "Плацкарт"

# New row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = "dsada"
$ws.Range("D6").Value = "dsadasda"
$ws.Range("E6").Value = 1000
$ws.Range("F6").Value = "0:0"
$ws.Range("G6").Value = "Плацкарт"

# New row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = "dsad"
$ws.Range("D7").Value = "asdsadas"
$ws.Range("E7").Value = 1000
$ws.Range("F7").Value = "0:0"
$ws.Range("G7").Value = "Плацкарт"

# New row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = "dsadasd"
$ws.Range("D8").Value = "sadasdas"
$ws.Range("E8").Value = 2001
$ws.Range("F8").Value = "0:0"
$ws.Range("G8").Value = "Плацкарт"

# Apply style to new A6:A8 cells (same style as A2:A5, style index 1 - centered bold border)
$ws.Range("A2").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)  # xlPasteFormats
